# Updates cryptos list cell values (Coin/Link/Price/Volume(1h)) to match the
# latest scrape, per commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.563.08'
$ws.Range('E2').Value = '  -0.35%  '
$ws.Range('D3').Value = '1.637.29'
$ws.Range('E3').Value = '  +0.39%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Value = "'213.77"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.66%  '
$ws.Range('D6').Value = "'0.505"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.24%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = "'0.251"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.25%  '
$ws.Range('D9').Value = "'0.0626"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.52%  '
$ws.Range('D10').Value = "'18.87"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.75%  '
$ws.Range('E11').Value = '  +0.58%  '
$ws.Range('D12').Value = '1.857.04'
$ws.Range('E12').Value = '  -0.08%  '
$ws.Range('D13').Value = '1.656.60'
$ws.Range('E13').Value = '  +1.16%  '
$ws.Range('D14').Value = "'4.16"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.02%  '
$ws.Range('D15').Value = "'0.526"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.06%  '
$ws.Range('D16').Value = "'65.41"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.89%  '
$ws.Range('D17').Value = '26.599.68'
$ws.Range('E17').Value = '  -0.12%  '
$ws.Range('D18').Value = '0.0₃0745'
$ws.Range('E18').Value = '  +0.97%  '
$ws.Range('D19').Value = "'216.10"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.34%  '
$ws.Range('E20').Value = '  +0.21%  '
$ws.Range('D21').Value = "'4.31"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.53%  '
$ws.Range('D22').Value = "'6.28"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.77%  '
$ws.Range('D23').Value = "'9.36"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.47%  '
$ws.Range('D24').Value = "'2.20"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +13.46%  '
$ws.Range('D25').Value = "'147.30"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.36%  '
$ws.Range('E26').Value = '  +0.09%  '
$ws.Range('E27').Value = '  +0.24%  '
$ws.Range('D28').Value = "'6.94"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.06%  '
$ws.Range('D29').Value = "'15.69"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.35%  '
$ws.Range('D30').Value = "'0.0515"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.06%  '
$ws.Range('E31').Value = '  -0.30%  '
$ws.Range('E32').Value = '  +4.26%  '
$ws.Range('D33').Value = "'2.99"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.53%  '
$ws.Range('D34').Value = '1.258.36'
$ws.Range('E34').Value = '  +7.78%  '
$ws.Range('D35').Value = "'1.51"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.07%  '
$ws.Range('D36').Value = "'2.39"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.41%  '
$ws.Range('E37').Value = '  +4.82%  '
$ws.Range('D38').Value = "'0.512"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.83%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').Value = "'0.800"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.90%  '
$ws.Range('B40').Value = 'PaxDollar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D40').Value = "'1.00"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.16%  '
$ws.Range('D41').Value = "'2.28"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.60%  '
$ws.Range('D42').Value = "'0.799"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.83%  '
$ws.Range('D43').Value = "'5.35"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.36%  '
$ws.Range('D44').Value = '1.765.90'
$ws.Range('E44').Value = '  -0.13%  '
$ws.Range('D45').Value = "'93.41"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.57%  '
$ws.Range('D46').Value = "'1.60"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.57%  '
$ws.Range('D47').Value = "'55.19"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.06%  '
$ws.Range('E48').Value = '  -1.93%  '
$ws.Range('D49').Value = "'0.0511"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.23%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = "'7.53"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.28%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').Value = "'0.408"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.39%  '
